$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header updates
# B3 holds a date written as DD/MM/YYYY free text. Mark the cell as Text
# first so the COM layer keeps it as a literal string instead of silently
# converting it into a date serial number.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "03/09/2025"
$ws.Range("B5").Value = "Formation de relais"

# Updated list of names (rows 8-30), numbers in column A (1-23) stay unchanged
$ws.Range("B8").Value  = "SEGLA MODESTE GBETOGNON"
$ws.Range("B9").Value  = "DEGUENON CINSOU JUDES CHARLES"
$ws.Range("B10").Value = "OUENSA GEO FROY MWILFRIED"
$ws.Range("B11").Value = "AGOI MARTIN"
$ws.Range("B12").Value = "AMOUSSOU PARFAIT"
$ws.Range("B13").Value = "ADEDE KOCOU ABRAHAM"
$ws.Range("B14").Value = "HOUNTONNAGNON CAKPO LUCIEN"
$ws.Range("B15").Value = "TCHOGNINOU MATHIAS"
$ws.Range("B16").Value = "ZINHOUEHOU DOROTHE"
$ws.Range("B17").Value = "FANOU SENOUMATE GEDEON"
$ws.Range("B18").Value = "EKEHOUNDE ECUROSSE MICHELINE"
$ws.Range("B19").Value = "DAGBEGNON BONAVENTURE"
$ws.Range("B20").Value = "SOHOUNDJO MEREMIE DENANHOUEA"
$ws.Range("B21").Value = "ADOKO LEA SONAGNON"
$ws.Range("B22").Value = "DOSSA ASYLVESTRE"
$ws.Range("B23").Value = "TOBOSSI REFI"
$ws.Range("B24").Value = "DOVONON CLEMENTINE"
$ws.Range("B25").Value = "GANSE SEWANOU ALBERTINE"
$ws.Range("B26").Value = "KOUADOUA JULIETTE AYABA M"
$ws.Range("B27").Value = "LAWANI KADER"
$ws.Range("B28").Value = "AGNIZO VICTORINE"
$ws.Range("B29").Value = "DOSSOUHOUI AHODEGNON ARMAND"
$ws.Range("B30").Value = "TCHOBO FAABO JANVIENNE"

# Remove the rows for person #24 through #44 (rows 31-51). This shifts the
# TOTAL row (was 52) up to row 31, and the signature rows (were 54/55) up
# to rows 33/34, matching the target layout exactly.
$ws.Rows("31:51").Delete()

# Update the signature block (now at rows 33-34)
$ws.Range("B34").Value = "AHODEKON Maxiès"
$ws.Range("E34").Value = "."
